# Edit reference names of tumor vessel parameters so that they contain
# mice strain and tumor cell types.

$wb = $excel.ActiveWorkbook

$wsSize    = $wb.Worksheets.Item("Vessel size (tumor)")
$wsDensity = $wb.Worksheets.Item("Vessel density (tumor)")

# --- "Vessel size (tumor)": rows unique to this sheet first --------------
$wsSize.Range("A2").Value = "Mesquita et al., 2012 (C3H mice & fibrosarcoma)"
$wsSize.Range("A3").Value = "Mesquita et al., 2012 (Nude mice & fibrosarcoma)"

# --- "Vessel density (tumor)": full top-to-bottom pass --------------------
$wsDensity.Range("A2").Value  = "Tufto & Rofstad, 1999 (Balb/c nu/nu & D-12 cell)"
$wsDensity.Range("A3").Value  = "Tufto & Rofstad, 1999 (Balb/c nu/nu & R-18 cell)"
$wsDensity.Range("A4").Value  = "Tufto & Rofstad, 1999 (Balb/c nu/nu & U-25 cell)"
$wsDensity.Range("A5").Value  = "Fernandez-Rodrigues et al., 2016 (C57Bl/6 & B16F1)"
$wsDensity.Range("A6").Value  = "Kostourou et al., 2013 (C57BL6/129 & B16F0 cell)"
$wsDensity.Range("A7").Value  = "Kostourou et al., 2013 (C57BL6/129 & CMT19T cell)"
$wsDensity.Range("A8").Value  = "Jones et al., 2013 (C57BL6 & LLC cell)"
$wsDensity.Range("A9").Value  = "Goel et al., 2013 (Nude & 4T1 primary tumor)"
$wsDensity.Range("A10").Value = "Goel et al., 2013 (C57BL6/J & E0771 tumor)"

$wsDensity.Columns.Item(1).ColumnWidth = 55.33203125
$wsDensity.Range("C17").Select()
$wsDensity.Activate()

# --- back to "Vessel size (tumor)": remaining rows (reuse existing text) --
$wsSize.Range("A4").Value = "Goel et al., 2013 (Nude & 4T1 primary tumor)"
$wsSize.Range("A5").Value = "Goel et al., 2013 (C57BL6/J & E0771 tumor)"
$wsSize.Range("A6").Value = "Tufto & Rofstad, 1999 (Balb/c nu/nu & D-12 cell)"
$wsSize.Range("A7").Value = "Tufto & Rofstad, 1999 (Balb/c nu/nu & R-18 cell)"
$wsSize.Range("A8").Value = "Tufto & Rofstad, 1999 (Balb/c nu/nu & U-25 cell)"

$wsSize.Columns.Item(1).ColumnWidth = 45.1640625
$wsSize.Range("A6").Select()

# --- Workbook window view --------------------------------------------------
$excel.Width = 22400
$excel.Height = 22900
$excel.Left = 22400
$excel.Top = 500
